# The edit moves the two topmost data rows (rows 2 and 3, directly below
# the header row) down to the bottom of the data range, while all the
# other data rows (originally rows 4-36) shift up by two rows to fill the
# gap. Row 36 stays the last row (same used range, just reordered).
#
# Implementation: copy rows 2:3 to a scratch area right after the last
# existing row (37:38), then delete the original rows 2:3. Deleting those
# rows shifts everything below (including the just-copied scratch rows)
# up by two, which leaves the copied rows sitting at 35:36 - exactly
# where they need to be - while rows 4-36 land on 2-34.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$lastCol = $ws.UsedRange.Columns.Count

$srcRange = $ws.Range($ws.Cells.Item(2, 1), $ws.Cells.Item(3, $lastCol))
$dstRange = $ws.Range($ws.Cells.Item($lastRow + 1, 1), $ws.Cells.Item($lastRow + 2, $lastCol))

$srcRange.Copy($dstRange) | Out-Null

$ws.Rows("2:3").Delete() | Out-Null
